# Update cryptocurrency price/volume snapshot data to match the latest
# scrape (commit: "Updated symbol list ... with GitHub Actions").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '244.50'

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '26.50'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '4.53%'

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.127'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '0.16%'

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05588'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '0.29%'

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.466'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-0.48%'

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8179'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '0.08%'

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8347'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-1.22%'

# Row 9
$ws.Range("B9").Value = 'WazirX'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1335'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-0.76%'

# Row 10
$ws.Range("B10").Value = 'MandalaExchangeToken'
$ws.Range("C10").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06985'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '0.39%'

# Row 11
$ws.Range("B11").Value = 'BitrueCoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.02888'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '0.09%'

# Row 12
$ws.Range("B12").Value = 'BitMartToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09387'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '0.11%'

# Row 13
$ws.Range("B13").Value = 'BitForexToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.001507'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.45%'

# Row 14
$ws.Range("B14").Value = 'TigerCash'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.006150'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '0.31%'

# Row 15
$ws.Range("B15").Value = 'LEO'
$ws.Range("C15").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.649'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '4.25%'

# Row 16
$ws.Range("B16").Value = 'GateToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.036'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '0.53%'

# Row 17
$ws.Range("B17").Value = 'BTSEToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.183'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '5.81%'

# Row 18
$ws.Range("B18").Value = 'One'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0005951'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-0.85%'

# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-2.11%'

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.03111'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-2.03%'

# Row 21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-2.24%'

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.744'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '0.02%'

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04670'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-1.10%'

# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-0.09%'

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001245'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-0.31%'

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004501'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-2.96%'

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.00009602'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-1.06%'

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '0.36%'

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03642'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-0.77%'

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006191'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '83.16%'

# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-22.66%'

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002400'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-9.46%'

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008874'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '6.56%'

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005360'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '1.24%'

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-0.02%'

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-4.02%'

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002335'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '10.05%'

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002100'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-0.02%'

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0002000'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.02%'
